$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 7 new blank rows right before the existing data block (old rows 782-788),
#     pushing that block (and everything below it) down by 7 rows. ---
$ws.Rows("782:788").Insert(-4121)

# --- Step 2: fill the 7 newly inserted rows (782-788) with the refreshed data values ---
$ws.Range("A782").Value = 44246
$ws.Range("B782").Value = 23323.98

$ws.Range("A783").Value = 44247
$ws.Range("A784").Value = 44248

# Give B783:B784 the "--" placeholder style (same as used elsewhere, e.g. B776:B777)
$ws.Range("B776").Copy()
$ws.Range("B783:B784").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B783").Value = "--"
$ws.Range("B784").Value = "--"

$ws.Range("A785").Value = 44249
$ws.Range("B785").Value = 23080.93

$ws.Range("A786").Value = 44250
$ws.Range("B786").Value = 22709.56

$ws.Range("A787").Value = 44251
$ws.Range("B787").Value = 23101.81

$ws.Range("A788").Value = 44252
$ws.Range("B788").Value = 23110.73

# --- Step 3: fill the new trailing data row (795) which is the first blank row after the
#     (now shifted) original data block ending at row 794. ---
$ws.Range("A795").Value = 44256
$ws.Range("B795").Value = 23403.58

# Reflect the new "latest entry" as the active selection (matches the author's last edit spot)
$ws.Range("B795").Select()

# --- Step 4: update the workbook-level defined name range to cover the new data extent ---
$n = $wb.Names.Item("IGPA")
$n.RefersTo = "=IGPA!`$A`$1:`$B`$788"

Write-Host "done"
